$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Telefone" column at D (inherits header formatting
#     from the neighbouring column, same as a real Excel column insert). ---
$ws.Columns("D").Insert()
$ws.Range("D1").Value = "Telefone"

# --- Reorder rows 4-6: "magic gril" moves up to row 4, pushing
#     "melanina" down to row 5 and "centauro west" down to row 6. ---
$ws.Range("A4").Value = "magic gril "
$ws.Range("B4").Value = "Rua Aricuri, 649"
$ws.Range("C4").Value = ""
$ws.Range("C4").Style = "Normal"

$ws.Range("A5").Value = "melanina "
$ws.Range("B5").Value = "R Soldado Felisbino dos Santos, 97 "
$ws.Range("C5").Value = "sobreloja "

$ws.Range("A6").Value = "centauro west "
$ws.Range("B6").Value = " Estr. do Mendanha, 555"
